$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$A2 = @'
Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),
                ('model',
                 MLPClassifier(activation='tanh', alpha=1e-05,
                               hidden_layer_sizes=(5, 10, 5),
                               learning_rate_init=0.01, max_iter=1000,
                               random_state=42, solver='sgd'))])
'@
$ws.Range("A2").Value = $A2

$ws.Range("B2").Value = 0.6476190476190476
$C2 = @'
{'scaler': None, 'model__solver': 'sgd', 'model__learning_rate_init': 0.01, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 1e-05, 'model__activation': 'tanh'}
'@
$ws.Range("C2").Value = $C2
$ws.Range("D2").Value = 0.5000000000000001
$ws.Range("E2").Value = "[1 0 0 1 0 0 1 1 0 1 0 0]"
$ws.Range("F2").Value = "[1 1 1 0 1 1 1 1 1 1 1 1]"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.7399523809523809
$ws.Range("I2").Value = 0.03062183231262731
$ws.Range("J2").Value = 0.5441904761904762
$ws.Range("K2").Value = 0.07866228171909134

# Row 3
$A3 = @'
Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),
                ('model',
                 MLPClassifier(activation='tanh', alpha=0.01,
                               hidden_layer_sizes=(10,),
                               learning_rate_init=1e-05, max_iter=1000,
                               random_state=42, solver='lbfgs'))])
'@
$ws.Range("A3").Value = $A3

$ws.Range("B3").Value = 0.638095238095238
$C3 = @'
{'scaler': None, 'model__solver': 'lbfgs', 'model__learning_rate_init': 1e-05, 'model__hidden_layer_sizes': (10,), 'model__alpha': 0.01, 'model__activation': 'tanh'}
'@
$ws.Range("C3").Value = $C3
$ws.Range("D3").Value = 0.625
$ws.Range("E3").Value = "[1 0 1 0 0 0 0 1 1 0 1 1]"
$ws.Range("F3").Value = "[1 1 1 1 1 0 1 1 0 1 1 1]"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.7247380952380954
$ws.Range("I3").Value = 0.03259983005491219
$ws.Range("J3").Value = 0.5303809523809524
$ws.Range("K3").Value = 0.06421710738702992

# Row 4
$A4 = @'
Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),
                ('model',
                 MLPClassifier(activation='tanh', alpha=0.01,
                               hidden_layer_sizes=(5, 10, 5),
                               learning_rate_init=0.01, max_iter=1000,
                               random_state=42, solver='sgd'))])
'@
$ws.Range("A4").Value = $A4

$ws.Range("B4").Value = 0.6190476190476191
$C4 = @'
{'scaler': None, 'model__solver': 'sgd', 'model__learning_rate_init': 0.01, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 0.01, 'model__activation': 'tanh'}
'@
$ws.Range("C4").Value = $C4
$ws.Range("D4").Value = 0.8
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 1 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.7310238095238095
$ws.Range("I4").Value = 0.02514326533843337
$ws.Range("J4").Value = 0.5281904761904762
$ws.Range("K4").Value = 0.06704702718624295

